$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.034.27'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('D3').Value = '2.633.49'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.38'
$ws.Range('E5').Value = '  +6.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.50'
$ws.Range('E6').Value = '  +2.93%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.602'
$ws.Range('E8').Value = '  +6.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.87'
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.103'
$ws.Range('E10').Value = '  +3.04%  '
$ws.Range('E11').Value = '  +5.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.340'
$ws.Range('E12').Value = '  +1.93%  '
$ws.Range('D13').Value = '3.097.47'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Value = '60.100.65'
$ws.Range('E14').Value = '  +2.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.41'
$ws.Range('E15').Value = '  +3.23%  '
$ws.Range('D16').Value = '2.644.44'
$ws.Range('E16').Value = '  +1.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.53'
$ws.Range('E18').Value = '  +3.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '342.10'
$ws.Range('E19').Value = '  +2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.32'
$ws.Range('E20').Value = '  +2.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.30'
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.94'
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.433'
$ws.Range('E24').Value = '  +5.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.165'
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.992'
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.33'
$ws.Range('E27').Value = '  +2.92%  '
$ws.Range('D28').Value = '0.0₃0775'
$ws.Range('E28').Value = '  +6.91%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.70'
$ws.Range('E30').Value = '  +4.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.14'
$ws.Range('E31').Value = '  +6.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '156.92'
$ws.Range('E32').Value = '  +4.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.13'
$ws.Range('E33').Value = '  +2.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.07'
$ws.Range('E34').Value = '  +4.89%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.18'
$ws.Range('E35').Value = '  +7.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.911'
$ws.Range('E36').Value = '  +11.72%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.909'
$ws.Range('E37').Value = '  +12.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.43'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.51'
$ws.Range('E39').Value = '  +6.28%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '298.83'
$ws.Range('E40').Value = '  +6.38%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.67'
$ws.Range('E41').Value = '  +3.41%  '
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.606'
$ws.Range('E43').Value = '  +2.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0548'
$ws.Range('E44').Value = '  +3.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0968'
$ws.Range('E45').Value = '  +3.64%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.40'
$ws.Range('E46').Value = '  +3.47%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.63'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0232'
$ws.Range('E48').Value = '  +4.00%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.73'
$ws.Range('E49').Value = '  +6.53%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.970.73'
$ws.Range('E50').Value = '  +2.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.46'
$ws.Range('E51').Value = '  +3.57%  '
